# Profimunka.pptx - apply the three small text tweaks described in the
# commit's OOXML diff:
#   1. Slide 2 ("Feladatunk"): "...keszitese volt" -> "...keszitese volt."
#   2. Slide 4 ("A Projekt celja"): "...szeretne megtudni" ->
#      "...szeretne megtudni:" (split into two runs at the same point
#      the author split them)
#   3. Slide 5 ("A Projekt kivitelezese"): merge the two runs
#      "Konzol " + "feluleten " back into a single run "Konzol feluleten "
#
# We edit via narrow TextRange.Characters(start, length) sub-ranges
# (found with IndexOf against the shape's full text) instead of
# clobbering the whole TextRange.Text, so the untouched runs/rPr in the
# paragraph are left completely alone and only the targeted run(s)
# change - matching the shape of the upstream diff as closely as this
# object model allows.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# 1) Slide 2, content placeholder: "volt" -> "volt."
# ---------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$body2 = $slide2.Shapes.Item(2)
$tr2 = $body2.TextFrame.TextRange
$text2 = $tr2.Text
$needle2 = "volt"
$pos2 = $text2.LastIndexOf($needle2)
if ($pos2 -ge 0) {
    $run2 = $tr2.Characters($pos2 + 1, $needle2.Length)
    $run2.Text = "volt."
}

# ---------------------------------------------------------------
# 2) Slide 4, content placeholder: split the last run so "megtudni"
#    becomes its own run reading "megtudni:"
# ---------------------------------------------------------------
$slide4 = $p.Slides.Item(4)
$body4 = $slide4.Shapes.Item(2)
$tr4 = $body4.TextFrame.TextRange
$text4 = $tr4.Text
$needle4 = "megtudni"
$pos4 = $text4.IndexOf($needle4)
if ($pos4 -ge 0) {
    $run4 = $tr4.Characters($pos4 + 1, $needle4.Length)
    $run4.Text = "megtudni:"
}

# ---------------------------------------------------------------
# 3) Slide 5, content placeholder: merge "Konzol " and "feluleten "
#    into the single run "Konzol feluleten "
# ---------------------------------------------------------------
$slide5 = $p.Slides.Item(5)
$body5 = $slide5.Shapes.Item(2)
$tr5 = $body5.TextFrame.TextRange
$merged5 = "Konzol felületen "
$run5 = $tr5.Characters(1, $merged5.Length)
$run5.Text = $merged5
